$p = $ppt.ActivePresentation

# --- 1) Retarget the three table styles (slides 14-16) from the old
#        "CF1EA01C..." style to the new "1B55764F..." style. ---
$newStyleId = "{1B55764F-9D28-4F71-9EE6-8258BB486069}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Swap the slide master's theme palette from the custom
#        "Integral" / "Red Violet" scheme to the default Office palette. ---
$master = $p.SlideMaster
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0x00000000
$scheme.Colors(2).RGB  = 0x00FFFFFF
$scheme.Colors(3).RGB  = 0x006A5444
$scheme.Colors(4).RGB  = 0x00E6E6E7
$scheme.Colors(5).RGB  = 0x00D59B5B
$scheme.Colors(6).RGB  = 0x00317DED
$scheme.Colors(7).RGB  = 0x00A5A5A5
$scheme.Colors(8).RGB  = 0x0000C0FF
$scheme.Colors(9).RGB  = 0x00C47244
$scheme.Colors(10).RGB = 0x0047AD70
$scheme.Colors(11).RGB = 0x00C16305
$scheme.Colors(12).RGB = 0x00724F95
